$d = $word.ActiveDocument

# --- Rebuild the body content (all paragraphs up to, but excluding, the
#     document's permanent final paragraph mark) via InsertXML so we get
#     precise control of the resulting run/paragraph structure -----------
$xmlBody =
  '<w:p><w:pPr><w:pStyle w:val="Titel"/><w:jc w:val="center"/><w:rPr/></w:pPr>' +
    '<w:r><w:rPr/><w:t>Ordbog</w:t></w:r></w:p>' +
  '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:rPr/></w:pPr>' +
    '<w:r><w:rPr/><w:t>Aftale</w:t></w:r></w:p>' +
  '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:rPr/></w:pPr>' +
    '<w:r><w:rPr/><w:tab/><w:t>Er et tilbud som er godkendt af bils&#230;lger/salgschef samt kunden.</w:t></w:r></w:p>' +
  '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:rPr/></w:pPr>' +
    '<w:r><w:rPr/><w:t>Tilbud</w:t></w:r></w:p>' +
  '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:rPr/></w:pPr>' +
    '<w:r><w:rPr/><w:tab/><w:t>Er et udkast til en aftale, som er godkendt af bils&#230;lger/salgschef.</w:t></w:r></w:p>' +
  '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:tabs><w:tab w:val="left" w:pos="6612" w:leader="none"/></w:tabs><w:rPr/></w:pPr>' +
    '<w:r><w:rPr/><w:t>FFS &#8211; Ferrari Financing System</w:t><w:tab/></w:r></w:p>'

$xmlPkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' + $xmlBody + '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'

$target = $d.Range(0, $d.Content.End - 1)
$target.InsertXML($xmlPkg)

# --- Section/page-setup tweaks (header & footer distance 708 -> 0) ------
$sec = $d.Sections(1)
$ps = $sec.PageSetup
$ps.HeaderDistance = 0
$ps.FooterDistance = 0
